# Update "想去人数" (column F) counts that changed between the two site
# generations (old commit -> 456a3b4).
#
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list mostly the
# same events, and sheet "演出" (sheet2) contributes one row that is also
# mirrored inside "全部类型". All three sheets need the corresponding F
# cells bumped to match.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$wsExhibit.Range("F3").Value  = 20370
$wsExhibit.Range("F12").Value = 43
$wsExhibit.Range("F15").Value = 18
$wsExhibit.Range("F17").Value = 201
$wsExhibit.Range("F19").Value = 453
$wsExhibit.Range("F26").Value = 1127
$wsExhibit.Range("F27").Value = 37
$wsExhibit.Range("F28").Value = 21
$wsExhibit.Range("F29").Value = 191
$wsExhibit.Range("F30").Value = 5219
$wsExhibit.Range("F31").Value = 575
$wsExhibit.Range("F32").Value = 85
$wsExhibit.Range("F33").Value = 4889
$wsExhibit.Range("F35").Value = 89
$wsExhibit.Range("F37").Value = 12744
$wsExhibit.Range("F39").Value = 94
$wsExhibit.Range("F43").Value = 385
$wsExhibit.Range("F44").Value = 4014
$wsExhibit.Range("F46").Value = 97

# --- Sheet "演出" ---
$wsShow.Range("F2").Value = 208

# --- Sheet "全部类型" ---
$wsAll.Range("F3").Value  = 20370
$wsAll.Range("F12").Value = 43
$wsAll.Range("F15").Value = 18
$wsAll.Range("F17").Value = 201
$wsAll.Range("F19").Value = 453
$wsAll.Range("F26").Value = 1127
$wsAll.Range("F27").Value = 37
$wsAll.Range("F28").Value = 21
$wsAll.Range("F29").Value = 191
$wsAll.Range("F30").Value = 208
$wsAll.Range("F31").Value = 5219
$wsAll.Range("F32").Value = 575
$wsAll.Range("F34").Value = 85
$wsAll.Range("F36").Value = 4889
$wsAll.Range("F38").Value = 89
$wsAll.Range("F40").Value = 12744
$wsAll.Range("F42").Value = 94
$wsAll.Range("F46").Value = 385
$wsAll.Range("F47").Value = 4014
$wsAll.Range("F49").Value = 97
